$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($rowIndex, $values) {
    $col = 1
    foreach ($v in $values) {
        $ws.Cells.Item($rowIndex, $col).Value = $v
        $col = $col + 1
    }
}

# Column layout: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric statistic columns.

# Row 2: ECs / Ntf5 / Ntrk1 / ECs
Set-RowValues 2 @("ECs", "Ntf5", "Ntrk1", "ECs", 1, 0.3333333333333333, 0.1277066666666667, 0.38312, 0.1657851361976475, 0.1974405699349423, 1, 0.3333333333333333, 0.0009183333333333334, 0.002755, 0.005182370534790543, 0.005409247358696911, 0.0001172772888888889, 0.0010554956, 0.0008591600049369252, 0.0010680048814202)

# Row 3: ECs / Ntf5 / Ntrk1 / FAPs
Set-RowValues 3 @("ECs", "Ntf5", "Ntrk1", "FAPs", 1, 0.3333333333333333, 0.1277066666666667, 0.38312, 0.1657851361976475, 0.1974405699349423, 2, 0.6666666666666666, 0.153988, 0.461964, 0.8689904253117887, 0.9070335923096405, 0.01966529418666666, 0.17698764768, 0.1440656960147665, 0.1790852294157536)

# Row 4: ECs / Ntf5 / Ntrk1 / MuSCs
Set-RowValues 4 @("ECs", "Ntf5", "Ntrk1", "MuSCs", 1, 0.3333333333333333, 0.1277066666666667, 0.38312, 0.1657851361976475, 0.1974405699349423, 1, 0.5, 0.022297, 0.044594, 0.1258272041534207, 0.08755716033166246, 0.002847475546666667, 0.01708485328, 0.02086028017794405, 0.01728733563776856)

# Row 5: MuSCs / Ntf5 / Ntrk1 / ECs
Set-RowValues 5 @("MuSCs", "Ntf5", "Ntrk1", "ECs", 1, 0.5, 0.370511, 0.741022, 0.4809867660085082, 0.3818850647690823, 1, 0.3333333333333333, 0.0009183333333333334, 0.002755, 0.005182370534790543, 0.005409247358696911, 0.0003402526016666667, 0.00204151561, 0.002492651643786687, 0.002065710777927958)

# Row 6: MuSCs / Ntf5 / Ntrk1 / FAPs
Set-RowValues 6 @("MuSCs", "Ntf5", "Ntrk1", "FAPs", 1, 0.5, 0.370511, 0.741022, 0.4809867660085082, 0.3818850647690823, 2, 0.6666666666666666, 0.153988, 0.461964, 0.8689904253117887, 0.9070335923096405, 0.05705424786799999, 0.342325487208, 0.4179728943630753, 0.3463825821469005)

# Row 7: MuSCs / Ntf5 / Ntrk1 / MuSCs
Set-RowValues 7 @("MuSCs", "Ntf5", "Ntrk1", "MuSCs", 1, 0.5, 0.370511, 0.741022, 0.4809867660085082, 0.3818850647690823, 1, 0.5, 0.022297, 0.044594, 0.1258272041534207, 0.08755716033166246, 0.008261283766999999, 0.033045135068, 0.06052122000164618, 0.03343677184425384)

# Row 8: Resolving-Mac / Ntf5 / Ntrk1 / ECs  (new row)
Set-RowValues 8 @("Resolving-Mac", "Ntf5", "Ntrk1", "ECs", 1, 0.3333333333333333, 0.2720966666666667, 0.81629, 0.3532280977938443, 0.4206743652959753, 1, 0.3333333333333333, 0.0009183333333333334, 0.002755, 0.005182370534790543, 0.005409247358696911, 0.0002498754388888889, 0.00224887895, 0.001830558886066931, 0.002275531699348754)

# Row 9: Resolving-Mac / Ntf5 / Ntrk1 / FAPs  (new row)
Set-RowValues 9 @("Resolving-Mac", "Ntf5", "Ntrk1", "FAPs", 1, 0.3333333333333333, 0.2720966666666667, 0.81629, 0.3532280977938443, 0.4206743652959753, 2, 0.6666666666666666, 0.153988, 0.461964, 0.8689904253117887, 0.9070335923096405, 0.04189962150666666, 0.37709659356, 0.3069518349339469, 0.3815657807469865)

# Row 10: Resolving-Mac / Ntf5 / Ntrk1 / MuSCs  (new row)
Set-RowValues 10 @("Resolving-Mac", "Ntf5", "Ntrk1", "MuSCs", 1, 0.3333333333333333, 0.2720966666666667, 0.81629, 0.3532280977938443, 0.4206743652959753, 1, 0.5, 0.022297, 0.044594, 0.1258272041534207, 0.08755716033166246, 0.006066939376666667, 0.03640163626, 0.04444570397383052, 0.03683305284964006)
